$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gnai2"
$ws.Cells.Item(2,3).Value = "F2r"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 203.7816646666667
$ws.Cells.Item(2,8).Value = 611.344994
$ws.Cells.Item(2,9).Value = 0.6667327591988204
$ws.Cells.Item(2,10).Value = 0.6667327591988205
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 8.038446
$ws.Cells.Item(2,14).Value = 24.115338
$ws.Cells.Item(2,15).Value = 0.1223314651384763
$ws.Cells.Item(2,16).Value = 0.1223314651384763
$ws.Cells.Item(2,17).Value = 1638.087907213108
$ws.Cells.Item(2,18).Value = 14742.79116491797
$ws.Cells.Item(2,19).Value = 0.0815623952886106
$ws.Cells.Item(2,20).Value = 0.08156239528861062

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gnai2"
$ws.Cells.Item(3,3).Value = "F2r"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 203.7816646666667
$ws.Cells.Item(3,8).Value = 611.344994
$ws.Cells.Item(3,9).Value = 0.6667327591988204
$ws.Cells.Item(3,10).Value = 0.6667327591988205
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 44.50790266666667
$ws.Cells.Item(3,14).Value = 133.523708
$ws.Cells.Item(3,15).Value = 0.6773345175739228
$ws.Cells.Item(3,16).Value = 0.6773345175739228
$ws.Cells.Item(3,17).Value = 9069.894496235305
$ws.Cells.Item(3,18).Value = 81629.05046611775
$ws.Cells.Item(3,19).Value = 0.4516011118026634
$ws.Cells.Item(3,20).Value = 0.4516011118026635

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Gnai2"
$ws.Cells.Item(4,3).Value = "F2r"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 203.7816646666667
$ws.Cells.Item(4,8).Value = 611.344994
$ws.Cells.Item(4,9).Value = 0.6667327591988204
$ws.Cells.Item(4,10).Value = 0.6667327591988205
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 13.16402266666667
$ws.Cells.Item(4,14).Value = 39.492068
$ws.Cells.Item(4,15).Value = 0.2003340172876008
$ws.Cells.Item(4,16).Value = 0.2003340172876008
$ws.Cells.Item(4,17).Value = 2682.586452723066
$ws.Cells.Item(4,18).Value = 24143.27807450759
$ws.Cells.Item(4,19).Value = 0.1335692521075462
$ws.Cells.Item(4,20).Value = 0.1335692521075463

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gnai2"
$ws.Cells.Item(5,3).Value = "F2r"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 63.14058933333333
$ws.Cells.Item(5,8).Value = 189.421768
$ws.Cells.Item(5,9).Value = 0.2065833519051582
$ws.Cells.Item(5,10).Value = 0.2065833519051582
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 8.038446
$ws.Cells.Item(5,14).Value = 24.115338
$ws.Cells.Item(5,15).Value = 0.1223314651384763
$ws.Cells.Item(5,16).Value = 0.1223314651384763
$ws.Cells.Item(5,17).Value = 507.552217764176
$ws.Cells.Item(5,18).Value = 4567.969959877584
$ws.Cells.Item(5,19).Value = 0.02527164411177544
$ws.Cells.Item(5,20).Value = 0.02527164411177544

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Gnai2"
$ws.Cells.Item(6,3).Value = "F2r"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 63.14058933333333
$ws.Cells.Item(6,8).Value = 189.421768
$ws.Cells.Item(6,9).Value = 0.2065833519051582
$ws.Cells.Item(6,10).Value = 0.2065833519051582
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 44.50790266666667
$ws.Cells.Item(6,14).Value = 133.523708
$ws.Cells.Item(6,15).Value = 0.6773345175739228
$ws.Cells.Item(6,16).Value = 0.6773345175739228
$ws.Cells.Item(6,17).Value = 2810.255204363972
$ws.Cells.Item(6,18).Value = 25292.29683927574
$ws.Cells.Item(6,19).Value = 0.1399260350014842
$ws.Cells.Item(6,20).Value = 0.1399260350014843

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Gnai2"
$ws.Cells.Item(7,3).Value = "F2r"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 63.14058933333333
$ws.Cells.Item(7,8).Value = 189.421768
$ws.Cells.Item(7,9).Value = 0.2065833519051582
$ws.Cells.Item(7,10).Value = 0.2065833519051582
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 13.16402266666667
$ws.Cells.Item(7,14).Value = 39.492068
$ws.Cells.Item(7,15).Value = 0.2003340172876008
$ws.Cells.Item(7,16).Value = 0.2003340172876008
$ws.Cells.Item(7,17).Value = 831.1841491706916
$ws.Cells.Item(7,18).Value = 7480.657342536224
$ws.Cells.Item(7,19).Value = 0.04138567279189847
$ws.Cells.Item(7,20).Value = 0.04138567279189848

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Gnai2"
$ws.Cells.Item(8,3).Value = "F2r"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 38.719942
$ws.Cells.Item(8,8).Value = 116.159826
$ws.Cells.Item(8,9).Value = 0.1266838888960214
$ws.Cells.Item(8,10).Value = 0.1266838888960214
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 8.038446
$ws.Cells.Item(8,14).Value = 24.115338
$ws.Cells.Item(8,15).Value = 0.1223314651384763
$ws.Cells.Item(8,16).Value = 0.1223314651384763
$ws.Cells.Item(8,17).Value = 311.248162890132
$ws.Cells.Item(8,18).Value = 2801.233466011188
$ws.Cells.Item(8,19).Value = 0.01549742573809025
$ws.Cells.Item(8,20).Value = 0.01549742573809025

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Gnai2"
$ws.Cells.Item(9,3).Value = "F2r"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 38.719942
$ws.Cells.Item(9,8).Value = 116.159826
$ws.Cells.Item(9,9).Value = 0.1266838888960214
$ws.Cells.Item(9,10).Value = 0.1266838888960214
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 44.50790266666667
$ws.Cells.Item(9,14).Value = 133.523708
$ws.Cells.Item(9,15).Value = 0.6773345175739228
$ws.Cells.Item(9,16).Value = 0.6773345175739228
$ws.Cells.Item(9,17).Value = 1723.343409794978
$ws.Cells.Item(9,18).Value = 15510.09068815481
$ws.Cells.Item(9,19).Value = 0.08580737076977508
$ws.Cells.Item(9,20).Value = 0.08580737076977511

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Gnai2"
$ws.Cells.Item(10,3).Value = "F2r"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 38.719942
$ws.Cells.Item(10,8).Value = 116.159826
$ws.Cells.Item(10,9).Value = 0.1266838888960214
$ws.Cells.Item(10,10).Value = 0.1266838888960214
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 13.16402266666667
$ws.Cells.Item(10,14).Value = 39.492068
$ws.Cells.Item(10,15).Value = 0.2003340172876008
$ws.Cells.Item(10,16).Value = 0.2003340172876008
$ws.Cells.Item(10,17).Value = 509.7101941400186
$ws.Cells.Item(10,18).Value = 4587.391747260168
$ws.Cells.Item(10,19).Value = 0.02537909238815604
$ws.Cells.Item(10,20).Value = 0.02537909238815605
